$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header values for new columns P (16th) and Q (17th) in row 1,
# matching the existing style used by the rest of row 1 (column O1, s="1").
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# For data rows 2-25: swap I/K and M/O column values, and fill in new P/Q columns.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2    # I column: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K column: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M column: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O column: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P column: new, value 2
    $ws.Cells.Item($r, 17).Value = 2   # Q column: new, value 2
}
